$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.143.66'
$ws.Range("E2").Value = '  -0.63%  '
$ws.Range("D3").Value = '3.156.35'
$ws.Range("E3").Value = '  -0.49%  '
$ws.Range("E4").Value = '  +0.32%  '
$ws.Range("D5").Value = "'613.12"
$ws.Range("E5").Value = '  +2.32%  '
$ws.Range("D6").Value = "'147.37"
$ws.Range("E6").Value = '  -2.85%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").Value = '3.151.78'
$ws.Range("E8").Value = '  -0.54%  '
$ws.Range("E9").Value = '  -0.67%  '
$ws.Range("E10").Value = '  -0.88%  '
$ws.Range("E11").Value = '  -2.11%  '
$ws.Range("E12").Value = '  -0.69%  '
$ws.Range("D13").Value = "'0.0000258"
$ws.Range("E13").Value = '  -0.38%  '
$ws.Range("D14").Value = "'35.53"
$ws.Range("E14").Value = '  -3.73%  '
$ws.Range("D15").Value = '3.677.53'
$ws.Range("E15").Value = '  -0.21%  '
$ws.Range("E16").Value = '  +2.94%  '
$ws.Range("D17").Value = '64.144.03'
$ws.Range("E17").Value = '  -0.58%  '
$ws.Range("D18").Value = '3.155.57'
$ws.Range("E18").Value = '  -0.59%  '
$ws.Range("D19").Value = "'6.89"
$ws.Range("E19").Value = '  -2.02%  '
$ws.Range("D20").Value = "'477.46"
$ws.Range("E20").Value = '  -0.61%  '
$ws.Range("D21").Value = "'14.68"
$ws.Range("E21").Value = '  -0.72%  '
$ws.Range("D22").Value = "'8.06"
$ws.Range("E22").Value = '  +4.04%  '
$ws.Range("D23").Value = "'0.715"
$ws.Range("E23").Value = '  -0.06%  '
$ws.Range("D24").Value = "'13.71"
$ws.Range("E24").Value = '  -1.22%  '
$ws.Range("D25").Value = "'83.65"
$ws.Range("E25").Value = '  -0.97%  '
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("D27").Value = "'2.83"
$ws.Range("E27").Value = '  -2.83%  '
$ws.Range("D28").Value = "'8.50"
$ws.Range("E28").Value = '  -1.47%  '
$ws.Range("D29").Value = "'7.09"
$ws.Range("E29").Value = '  +1.27%  '
$ws.Range("D30").Value = "'0.119"
$ws.Range("E30").Value = '  -1.54%  '
$ws.Range("D31").Value = "'2.11"
$ws.Range("E31").Value = '  -6.78%  '
$ws.Range("E32").Value = '  +0.25%  '
$ws.Range("D33").Value = "'2.70"
$ws.Range("E33").Value = '  -0.80%  '
$ws.Range("D34").Value = "'26.32"
$ws.Range("E34").Value = '  -2.06%  '
$ws.Range("E35").Value = '  +1.66%  '
$ws.Range("D36").Value = '0.0₃0786'
$ws.Range("E36").Value = '  +6.41%  '
$ws.Range("D37").Value = "'6.01"
$ws.Range("E37").Value = '  -1.89%  '
$ws.Range("D38").Value = "'52.94"
$ws.Range("E38").Value = '  -2.82%  '
$ws.Range("D39").Value = "'3.15"
$ws.Range("E39").Value = '  -2.40%  '
$ws.Range("D40").Value = "'462.15"
$ws.Range("E40").Value = '  +1.14%  '
$ws.Range("D41").Value = "'0.0399"
$ws.Range("E41").Value = '  -0.78%  '
$ws.Range("E42").Value = '  -3.94%  '
$ws.Range("D43").Value = "'8.37"
$ws.Range("E43").Value = '  -1.52%  '
$ws.Range("D44").Value = '2.864.48'
$ws.Range("E44").Value = '  -0.07%  '
$ws.Range("D45").Value = "'0.268"
$ws.Range("E45").Value = '  -2.41%  '
$ws.Range("E46").Value = '  -4.82%  '
$ws.Range("D47").Value = "'2.45"
$ws.Range("E47").Value = '  +4.59%  '
$ws.Range("D48").Value = "'26.54"
$ws.Range("E48").Value = '  -2.60%  '
$ws.Range("E49").Value = '  -0.05%  '
$ws.Range("E50").Value = '  -1.59%  '
$ws.Range("D51").Value = "'119.37"
$ws.Range("E51").Value = '  -0.30%  '
